# Updated placeholder index page and updated the products table to
# dynamically generate rows.
#
# This adds a new "img/background.jpg" asset row to the end of the
# Content Asset Log table (Asset / Source / Reasoning / Date Retrieved).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Append a brand new row at the end of the table. Word mirrors the
# tcPr/cnfStyle pattern of the existing rows automatically.
$newRow = $t.Rows.Add()

$cell1 = $newRow.Cells.Item(1)
$cell2 = $newRow.Cells.Item(2)
$cell3 = $newRow.Cells.Item(3)
$cell4 = $newRow.Cells.Item(4)

# First cell: "img/background.jpg" - authored with the same spell-check
# markup (proofErr spellStart/spellEnd around "img") that Word produces
# when it flags "img" as a non-dictionary word.
$cell1.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>img</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/background.jpg</w:t></w:r></w:p>')

$cell2.Range.Text = "https://unsplash.com/photos/aerial-photography-of-mall-interior-gFjGZ2qRZOo"
$cell3.Range.Text = "Background image for the context of retail"
$cell4.Range.Text = "12/01/2025"
